$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")
